$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Sheet1"

# Header row
$newSheet.Cells.Item(1,1).Value = "Index"
$newSheet.Cells.Item(1,2).Value = "Quantity"
$newSheet.Cells.Item(1,3).Value = "Part Number"
$newSheet.Cells.Item(1,4).Value = "Manufacturer Part Number"
$newSheet.Cells.Item(1,5).Value = "Description"
$newSheet.Cells.Item(1,6).Value = "Customer Reference"
$newSheet.Cells.Item(1,7).Value = "Available"
$newSheet.Cells.Item(1,8).Value = "Backorder"
$newSheet.Cells.Item(1,9).Value = "Unit Price"
$newSheet.Cells.Item(1,10).Value = "Extended Price USD"

# Row 2
$newSheet.Cells.Item(2,1).Value = 1
$newSheet.Cells.Item(2,2).Value = 10
$newSheet.Cells.Item(2,3).Value = "563-CFS-0102CT-ND"
$newSheet.Cells.Item(2,4).Value = "CFS-0102TB"
$newSheet.Cells.Item(2,5).Value = "SWITCH DIP 1POS SPST 100MA 6V"
$newSheet.Cells.Item(2,7).Value = 10
$newSheet.Cells.Item(2,8).Value = 0
$newSheet.Cells.Item(2,9).Value = 0.774
$newSheet.Cells.Item(2,10).Value = 7.74

# Row 3
$newSheet.Cells.Item(3,1).Value = 2
$newSheet.Cells.Item(3,2).Value = 10
$newSheet.Cells.Item(3,3).Value = "490-CSTNE16M0V53L000R0CT-ND"
$newSheet.Cells.Item(3,4).Value = "CSTNE16M0V53L000R0"
$newSheet.Cells.Item(3,5).Value = "3.2X1.3MM 16.0MHZ CERAMIC RESONA"
$newSheet.Cells.Item(3,7).Value = 10
$newSheet.Cells.Item(3,8).Value = 0
$newSheet.Cells.Item(3,9).Value = 0.225
$newSheet.Cells.Item(3,10).Value = 2.25

# Row 4
$newSheet.Cells.Item(4,1).Value = 3
$newSheet.Cells.Item(4,2).Value = 10
$newSheet.Cells.Item(4,3).Value = "LT1117CST-5#TRPBFCT-ND"
$newSheet.Cells.Item(4,4).Value = "LT1117CST-5#TRPBF"
$newSheet.Cells.Item(4,5).Value = "IC REG LINEAR 5V 800MA SOT223-3"
$newSheet.Cells.Item(4,7).Value = 10
$newSheet.Cells.Item(4,8).Value = 0
$newSheet.Cells.Item(4,9).Value = 5.172
$newSheet.Cells.Item(4,10).Value = 51.72

# Row 5
$newSheet.Cells.Item(5,1).Value = 4
$newSheet.Cells.Item(5,2).Value = 50
$newSheet.Cells.Item(5,3).Value = "SAM12304-ND"
$newSheet.Cells.Item(5,4).Value = "TSW-103-14-F-S"
$newSheet.Cells.Item(5,5).Value = "CONN HEADER VERT 3POS 2.54MM"
$newSheet.Cells.Item(5,7).Value = 50
$newSheet.Cells.Item(5,8).Value = 0
$newSheet.Cells.Item(5,9).Value = 0.381
$newSheet.Cells.Item(5,10).Value = 19.05

# Row 6
$newSheet.Cells.Item(6,1).Value = 5
$newSheet.Cells.Item(6,2).Value = 30
$newSheet.Cells.Item(6,3).Value = "S9337-ND"
$newSheet.Cells.Item(6,4).Value = "QPC02SXGN-RC"
$newSheet.Cells.Item(6,5).Value = "CONN JUMPER SHORTING .100`" GOLD"
$newSheet.Cells.Item(6,7).Value = 30
$newSheet.Cells.Item(6,8).Value = 0
$newSheet.Cells.Item(6,9).Value = 0.0408
$newSheet.Cells.Item(6,10).Value = 1.22

# Row 7
$newSheet.Cells.Item(7,1).Value = 6
$newSheet.Cells.Item(7,2).Value = 30
$newSheet.Cells.Item(7,3).Value = "277-5744-ND"
$newSheet.Cells.Item(7,4).Value = 1751264
$newSheet.Cells.Item(7,5).Value = "TERM BLK 4POS SIDE ENT 3.5MM PCB"
$newSheet.Cells.Item(7,7).Value = 30
$newSheet.Cells.Item(7,8).Value = 0
$newSheet.Cells.Item(7,9).Value = 2.749
$newSheet.Cells.Item(7,10).Value = 82.47

# Row 8
$newSheet.Cells.Item(8,1).Value = 7
$newSheet.Cells.Item(8,2).Value = 13
$newSheet.Cells.Item(8,3).Value = "Z2929-ND"
$newSheet.Cells.Item(8,4).Value = "G5Q-14 DC5"
$newSheet.Cells.Item(8,5).Value = "RELAY GEN PURPOSE SPDT 10A 5V"
$newSheet.Cells.Item(8,7).Value = 13
$newSheet.Cells.Item(8,8).Value = 0
$newSheet.Cells.Item(8,9).Value = 1.556
$newSheet.Cells.Item(8,10).Value = 20.23

# Row 9
$newSheet.Cells.Item(9,1).Value = 8
$newSheet.Cells.Item(9,2).Value = 10
$newSheet.Cells.Item(9,3).Value = "SAM11096-ND"
$newSheet.Cells.Item(9,4).Value = "ESQ-110-13-L-D"
$newSheet.Cells.Item(9,5).Value = "CONN SOCKET 20POS 0.1 GOLD PCB"
$newSheet.Cells.Item(9,7).Value = 10
$newSheet.Cells.Item(9,8).Value = 0
$newSheet.Cells.Item(9,9).Value = 6.01
$newSheet.Cells.Item(9,10).Value = 60.1

# Row 10
$newSheet.Cells.Item(10,1).Value = 9
$newSheet.Cells.Item(10,2).Value = 5
$newSheet.Cells.Item(10,3).Value = "ESQ-120-13-T-D-ND"
$newSheet.Cells.Item(10,4).Value = "ESQ-120-13-T-D"
$newSheet.Cells.Item(10,5).Value = "CONN SOCKET 40POS 0.1 TIN PCB"
$newSheet.Cells.Item(10,7).Value = 5
$newSheet.Cells.Item(10,8).Value = 0
$newSheet.Cells.Item(10,9).Value = 7.64
$newSheet.Cells.Item(10,10).Value = 38.2

# Row 11
$newSheet.Cells.Item(11,1).Value = 10
$newSheet.Cells.Item(11,2).Value = 8
$newSheet.Cells.Item(11,3).Value = "SAM1204-12-ND"
$newSheet.Cells.Item(11,4).Value = "SSQ-112-03-T-D"
$newSheet.Cells.Item(11,5).Value = "CONN RCPT 24POS 0.1 TIN PCB"
$newSheet.Cells.Item(11,7).Value = 8
$newSheet.Cells.Item(11,8).Value = 0
$newSheet.Cells.Item(11,9).Value = 2.56
$newSheet.Cells.Item(11,10).Value = 20.48

# Row 12
$newSheet.Cells.Item(12,1).Value = 11
$newSheet.Cells.Item(12,2).Value = 7
$newSheet.Cells.Item(12,3).Value = "SAM11124-ND"
$newSheet.Cells.Item(12,4).Value = "ESW-105-12-L-D"
$newSheet.Cells.Item(12,5).Value = "CONN SOCKET 10POS 0.1 GOLD PCB"
$newSheet.Cells.Item(12,7).Value = 7
$newSheet.Cells.Item(12,8).Value = 0
$newSheet.Cells.Item(12,9).Value = 2.99
$newSheet.Cells.Item(12,10).Value = 20.93

# Row 13
$newSheet.Cells.Item(13,1).Value = 12
$newSheet.Cells.Item(13,2).Value = 7
$newSheet.Cells.Item(13,3).Value = "609-77313-127-10LF-ND"
$newSheet.Cells.Item(13,4).Value = "77313-127-10LF"
$newSheet.Cells.Item(13,5).Value = "CONN HEADER VERT 10POS 2.54MM"
$newSheet.Cells.Item(13,7).Value = 7
$newSheet.Cells.Item(13,8).Value = 0
$newSheet.Cells.Item(13,9).Value = 0.62
$newSheet.Cells.Item(13,10).Value = 4.34

# Row 14
$newSheet.Cells.Item(14,1).Value = 13
$newSheet.Cells.Item(14,2).Value = 7
$newSheet.Cells.Item(14,3).Value = "S2212EC-03-ND"
$newSheet.Cells.Item(14,4).Value = "PREC003DFAN-RC"
$newSheet.Cells.Item(14,5).Value = "CONN HEADER VERT 6POS 2.54MM"
$newSheet.Cells.Item(14,7).Value = 7
$newSheet.Cells.Item(14,8).Value = 0
$newSheet.Cells.Item(14,9).Value = 0.16
$newSheet.Cells.Item(14,10).Value = 1.12

# Row 15
$newSheet.Cells.Item(15,1).Value = 14
$newSheet.Cells.Item(15,2).Value = 20
$newSheet.Cells.Item(15,3).Value = "296-52901-ND"
$newSheet.Cells.Item(15,4).Value = "LM358DG4"
$newSheet.Cells.Item(15,5).Value = "IC OPAMP GP 2 CIRCUIT 8SOIC"
$newSheet.Cells.Item(15,7).Value = 20
$newSheet.Cells.Item(15,8).Value = 0
$newSheet.Cells.Item(15,9).Value = 0.879
$newSheet.Cells.Item(15,10).Value = 17.58

# Row 16
$newSheet.Cells.Item(16,1).Value = 15
$newSheet.Cells.Item(16,2).Value = 7
$newSheet.Cells.Item(16,3).Value = "BCS-103-T-D-TE-ND"
$newSheet.Cells.Item(16,4).Value = "BCS-103-T-D-TE"
$newSheet.Cells.Item(16,5).Value = "CONN RCPT 6POS 0.1 TIN PCB"
$newSheet.Cells.Item(16,7).Value = 7
$newSheet.Cells.Item(16,8).Value = 0
$newSheet.Cells.Item(16,9).Value = 1.84
$newSheet.Cells.Item(16,10).Value = 12.88

# Row 17
$newSheet.Cells.Item(17,1).Value = 16
$newSheet.Cells.Item(17,2).Value = 10
$newSheet.Cells.Item(17,3).Value = "1195-6406-ND"
$newSheet.Cells.Item(17,4).Value = 21348081380050
$newSheet.Cells.Item(17,5).Value = "CBL 3POS MALE TO FMALE 16.4'"
$newSheet.Cells.Item(17,7).Value = 10
$newSheet.Cells.Item(17,8).Value = 0
$newSheet.Cells.Item(17,9).Value = 14.283
$newSheet.Cells.Item(17,10).Value = 142.83

$newSheet.Range("L23").Select()
